$wb = $excel.ActiveWorkbook

# Insert a new "State" column into hotel_info (between Hotel_Name and City)
# and populate it with the hotel's state.
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Cells.Item(1,3).Value = "State"
$hotelInfo.Cells.Item(2,3).Value = "Louisiana"

# Reorder the sheet tabs so review_info comes before hotel_info.
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($hotelInfo)
